$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ---------------------------------------------------------------------------
# 1. New log entries in rows 12-15 (continuing the "Leistungsdokumentation"
#    table). Columns: A=Tag (date), B=Von, C=Bis, D=Leistung,
#    E=Vollendung der Arbeit, F=Zusammenarbeit mit
#    (cells are touched in the same left-to-right / row-by-row order the
#    author used, with E14 entered last, so new shared-string entries land
#    in the same order as the original edit)
# ---------------------------------------------------------------------------
$ws.Cells.Item(12, 1).Value2 = 43972
$ws.Cells.Item(12, 2).Value2 = "20 Uhr"
$ws.Cells.Item(12, 3).Value2 = "21 Uhr"
$ws.Cells.Item(12, 4).Value2 = "Exploding Kittens probespielen, Brainstorming zur Funktionsweise eines Bot"
$ws.Cells.Item(12, 5).Value2 = "Grobes Konzept für schweren Bot im Kopf"
$ws.Cells.Item(12, 6).Value2 = "/"

$ws.Cells.Item(13, 1).Value2 = 43973
$ws.Cells.Item(13, 2).Value2 = "12 Uhr"
$ws.Cells.Item(13, 3).Value2 = "13:15 Uhr"
$ws.Cells.Item(13, 4).Value2 = "Schriftliche Erarbeitung Konzept für leichten und schweren Bot"
$ws.Cells.Item(13, 5).Value2 = "Abläufe für leichten und schweren Bot besprochen, ausgearbeitet und schriftlich festgehalten"
$ws.Cells.Item(13, 6).Value2 = "Manuela"

$ws.Cells.Item(14, 1).Value2 = 43973
$ws.Cells.Item(14, 2).Value2 = "16:15 Uhr"
$ws.Cells.Item(14, 3).Value2 = "17 Uhr"
$ws.Cells.Item(14, 4).Value2 = "Diagrammerstellung begonnen"
$ws.Cells.Item(14, 6).Value2 = "/"

$ws.Cells.Item(15, 1).Value2 = 43974
$ws.Cells.Item(15, 2).Value2 = "9 Uhr"
$ws.Cells.Item(15, 3).Value2 = "9:15 Uhr"
$ws.Cells.Item(15, 4).Value2 = "Diagrammerstellung abgeschlossen"
$ws.Cells.Item(15, 5).Value2 = "Diagramm für schweren Bot fertig"
$ws.Cells.Item(15, 6).Value2 = "/"

$ws.Cells.Item(14, 5).Value2 = "Diagramm für schweren Bot vorläufig fertig"

# ---------------------------------------------------------------------------
# 2. Re-apply the table's standard per-column formatting (date format /
#    borders / left alignment) down through row 25 so the newly used rows
#    (12-15) and the still-empty trailing rows (16-25) look like the rest
#    of the table instead of the old placeholder styling.
# ---------------------------------------------------------------------------
$ws.Range("A6").Copy()
$ws.Range("A12:A25").PasteSpecial(-4122)

$ws.Range("B6").Copy()
$ws.Range("B12:B25").PasteSpecial(-4122)

$ws.Range("C6").Copy()
$ws.Range("C12:C25").PasteSpecial(-4122)

$ws.Range("D6").Copy()
$ws.Range("D12:D25").PasteSpecial(-4122)

$ws.Range("E9").Copy()
$ws.Range("E12:E25").PasteSpecial(-4122)

$ws.Range("F6").Copy()
$ws.Range("F12:F25").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# D18 ends up with the date-column styling (matches the author's edit,
# where the active/selected cell picked up the neighbouring column format).
$ws.Range("A6").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Move the selection to D18, matching the saved cursor position.
# ---------------------------------------------------------------------------
[void]$ws.Range("D18").Select()
